$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.103.40"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "1.677.57"
$ws.Range("E3").Value = "  +3.76%  "
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "216.44"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  -0.02%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.252"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +3.18%  "
$ws.Range("E9").Value = "  +1.99%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "20.23"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +5.59%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0888"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("D12").Value = "1.914.22"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").Value = "1.675.97"
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("E14").Value = "  +1.84%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "66.02"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("E16").Value = "  +3.01%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "239.45"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "27.122.03"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "0.0₃0738"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +4.43%  "
$ws.Range("E23").Value = "  +2.88%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "9.30"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("E27").Value = "  +0.74%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "16.01"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  -0.12%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.0498"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "1.480.23"
$ws.Range("E33").Value = "  -3.13%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "3.12"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").Value = "  +6.44%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.903"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +8.92%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.575"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E39").Value = "  +2.46%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "6.06"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("E41").Value = "  -0.09%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.990"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +8.82%  "
$ws.Range("E43").Value = "  +4.09%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "66.79"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +8.69%  "
$ws.Range("D45").Value = "1.823.35"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  +2.16%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "90.41"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("E50").Value = "  +1.23%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "7.67"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +2.70%  "
